$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" (C) column for existing rows 2-196: 45182 -> 45184 ---
$ws.Range("C2:C196").Value2 = 45184

# --- 2. Row 196 gains an explicit row height (ht="15" customHeight="1") ---
$ws.Rows.Item(196).RowHeight = 15

# --- 3. Append three new data rows (197-199) ---

# Row 197
$ws.Range("A197").Value2 = "A 42554-2023"
$ws.Range("B197").Value2 = 45181
$ws.Range("C197").Value2 = 45184
$ws.Range("D197").Value2 = "NORRBOTTENS LÄN"
$ws.Range("E197").Value2 = "HAPARANDA"
$ws.Range("G197").Value2 = 2.1
$ws.Range("H197:Q197").Value2 = 0

# Row 198
$ws.Range("A198").Value2 = "A 42602-2023"
$ws.Range("B198").Value2 = 45181
$ws.Range("C198").Value2 = 45184
$ws.Range("D198").Value2 = "NORRBOTTENS LÄN"
$ws.Range("E198").Value2 = "HAPARANDA"
$ws.Range("G198").Value2 = 6.7
$ws.Range("H198:Q198").Value2 = 0

# Row 199
$ws.Range("A199").Value2 = "A 42598-2023"
$ws.Range("B199").Value2 = 45181
$ws.Range("C199").Value2 = 45184
$ws.Range("D199").Value2 = "NORRBOTTENS LÄN"
$ws.Range("E199").Value2 = "HAPARANDA"
$ws.Range("G199").Value2 = 2.6
$ws.Range("H199:Q199").Value2 = 0

# --- Formatting for the new rows, matching the existing table style ---

# B/C columns use the date number format (style index 1 => numFmtId 165 "YYYY-MM-DD")
$ws.Range("B197:C199").NumberFormat = "YYYY-MM-DD"

# R column uses the wrap-text style (style index 2) and is left blank
$ws.Range("R197:R199").WrapText = $true

# Rows 197 & 198 get an explicit row height like row 196; row 199 keeps the default (no explicit height)
$ws.Rows.Item(197).RowHeight = 15
$ws.Rows.Item(198).RowHeight = 15
